# Commit: "Pdf file to explain csv files is added"
#
# Functional spreadsheet changes (per the OOXML diff):
#   - Lithium_Battery  Eta_ch  (D10): 0.9  -> 0.95
#   - Lithium_Battery  Eta_dis (D11): 0.9  -> 0.95
#   - PCM (Thermal_Storage) Eta_ch  (D35): 80 -> 0.8
#   - PCM (Thermal_Storage) Eta_dis (D36): 80 -> 0.8
#   - Active cell/selection moved to D46 with the view scrolled so row 26 is
#     at the top of the window.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Lithium_Battery efficiencies
$ws.Range("D10").Value = 0.95
$ws.Range("D11").Value = 0.95

# PCM (thermal storage) efficiencies corrected from percentages written as
# whole numbers (80) to the proper fractional form (0.8)
$ws.Range("D35").Value = 0.8
$ws.Range("D36").Value = 0.8

# Restore the view/selection state saved in the workbook
$ws.Activate()
$ws.Range("D46").Select()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
